$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 7: score (C7) changes from 328 to 333 (B7/D7/E7 stay the same)
$ws.Range("C7").Value = 333

# Remove row 8 entirely (id=7, name="alice 333", score=387, tag_id="tag_1")
$ws.Rows.Item(8).Delete()

# Update the selection to match the new state
$ws.Range("C7").Select()

# Match the new zoom level recorded in the sheet view
$excel.ActiveWindow.Zoom = 99
